# Apply updated cryptocurrency price/volume figures (2023-04-27 GitHub Actions refresh).
# Each entry is a cell reference ("D<row>" / "E<row>") plus its new display text.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2"; Value = '29.572.86' },
    @{ Cell = "E2"; Value = '  +4.27%  ' },
    @{ Cell = "D3"; Value = '1.916.74' },
    @{ Cell = "E3"; Value = '  +2.62%  ' },
    @{ Cell = "D4"; Value = '1.001' },
    @{ Cell = "E4"; Value = '  -0.05%  ' },
    @{ Cell = "D5"; Value = '334.36' },
    @{ Cell = "E5"; Value = '  +1.09%  ' },
    @{ Cell = "E6"; Value = '  -0.05%  ' },
    @{ Cell = "D7"; Value = '0.4685' },
    @{ Cell = "E7"; Value = '  +2.15%  ' },
    @{ Cell = "D8"; Value = '0.4111' },
    @{ Cell = "E8"; Value = '  +2.86%  ' },
    @{ Cell = "D9"; Value = '48.13' },
    @{ Cell = "E9"; Value = '  +1.11%  ' },
    @{ Cell = "D10"; Value = '0.08047' },
    @{ Cell = "E10"; Value = '  +2.68%  ' },
    @{ Cell = "D11"; Value = '1.015' },
    @{ Cell = "E11"; Value = '  +3.21%  ' },
    @{ Cell = "D12"; Value = '22.40' },
    @{ Cell = "E12"; Value = '  +5.51%  ' },
    @{ Cell = "D13"; Value = '1.931.66' },
    @{ Cell = "E13"; Value = '  +3.00%  ' },
    @{ Cell = "D14"; Value = '5.998' },
    @{ Cell = "E14"; Value = '  +2.87%  ' },
    @{ Cell = "D15"; Value = '7.185' },
    @{ Cell = "E15"; Value = '  +2.84%  ' },
    @{ Cell = "D16"; Value = '89.98' },
    @{ Cell = "E17"; Value = '  -0.09%  ' },
    @{ Cell = "E18"; Value = '  +1.74%  ' },
    @{ Cell = "D19"; Value = '0.06587' },
    @{ Cell = "E19"; Value = '  +0.88%  ' },
    @{ Cell = "D20"; Value = '17.82' },
    @{ Cell = "E20"; Value = '  +3.88%  ' },
    @{ Cell = "E21"; Value = '  -0.11%  ' },
    @{ Cell = "D22"; Value = '29.550.31' },
    @{ Cell = "E22"; Value = '  +4.28%  ' },
    @{ Cell = "D23"; Value = '5.587' },
    @{ Cell = "E23"; Value = '  +4.85%  ' },
    @{ Cell = "E24"; Value = '  +6.77%  ' },
    @{ Cell = "D25"; Value = '2.210' },
    @{ Cell = "E25"; Value = '  -1.81%  ' },
    @{ Cell = "D26"; Value = '2.153.02' },
    @{ Cell = "E26"; Value = '  +2.53%  ' },
    @{ Cell = "D27"; Value = '155.41' },
    @{ Cell = "E27"; Value = '  -1.28%  ' },
    @{ Cell = "D28"; Value = '19.91' },
    @{ Cell = "E28"; Value = '  +3.15%  ' },
    @{ Cell = "D29"; Value = '5.768' },
    @{ Cell = "E29"; Value = '  +9.24%  ' },
    @{ Cell = "D30"; Value = '2.144' },
    @{ Cell = "E30"; Value = '  +4.30%  ' },
    @{ Cell = "D31"; Value = '117.65' },
    @{ Cell = "E31"; Value = '  +0.48%  ' },
    @{ Cell = "D32"; Value = '1.071' },
    @{ Cell = "E32"; Value = '  +12.35%  ' },
    @{ Cell = "D33"; Value = '0.09471' },
    @{ Cell = "E33"; Value = '  +1.78%  ' },
    @{ Cell = "E34"; Value = '  +3.41%  ' },
    @{ Cell = "E35"; Value = '  -0.42%  ' },
    @{ Cell = "D36"; Value = '5.419' },
    @{ Cell = "E36"; Value = '  +3.85%  ' },
    @{ Cell = "E37"; Value = '  +1.85%  ' },
    @{ Cell = "D39"; Value = '8.417' },
    @{ Cell = "E39"; Value = '  +1.78%  ' },
    @{ Cell = "D40"; Value = '1.178' },
    @{ Cell = "E40"; Value = '  +1.44%  ' },
    @{ Cell = "D41"; Value = '0.5903' },
    @{ Cell = "E41"; Value = '  +2.83%  ' },
    @{ Cell = "D42"; Value = '0.1843' },
    @{ Cell = "E42"; Value = '  +2.20%  ' },
    @{ Cell = "D43"; Value = '10.24' },
    @{ Cell = "E43"; Value = '  +2.36%  ' },
    @{ Cell = "D44"; Value = '1.274' },
    @{ Cell = "E44"; Value = '  +1.11%  ' },
    @{ Cell = "D45"; Value = '2.359' },
    @{ Cell = "E45"; Value = '  +3.84%  ' },
    @{ Cell = "D46"; Value = '0.07511' },
    @{ Cell = "E46"; Value = '  +5.24%  ' },
    @{ Cell = "D47"; Value = '12.24' },
    @{ Cell = "E47"; Value = '  +3.60%  ' },
    @{ Cell = "D48"; Value = '0.5576' },
    @{ Cell = "E48"; Value = '  +3.00%  ' },
    @{ Cell = "E49"; Value = '  +2.88%  ' },
    @{ Cell = "D50"; Value = '113.51' },
    @{ Cell = "E50"; Value = '  +3.22%  ' },
    @{ Cell = "D51"; Value = '0.2989' },
    @{ Cell = "E51"; Value = '  +12.13%  ' }
)

foreach ($u in $updates) {
    $cellRef = $u.Cell
    $newValue = $u.Value
    $range = $ws.Range($cellRef)

    if ($newValue -match '^-?\d+(\.\d+)?$') {
        # Looks like a plain decimal number (e.g. "334.36") - without forcing a
        # text format Excel would silently convert it to a numeric value, which
        # would not match the site's display string (and loses the original
        # trailing-zero formatting). Force text, write it, then strip the
        # number-format override back off so the cell style is left untouched.
        $range.NumberFormat = "@"
        $range.Value = $newValue
        $range.ClearFormats()
    } else {
        # Already unambiguous text (percent strings, multi-dot "thousands" style
        # prices like "29.572.86") - plain assignment keeps it as text.
        $range.Value = $newValue
    }
}
